# Refresh the cryptos list (prices / 1h volume %) on the "cryptos" sheet,
# mirroring the periodic GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be force-written as
# text (NumberFormat "@") so Excel does not coerce the string into a binary
# double (which would lose the exact decimal text, e.g. "134.45" -> "134.44999999999999").
# ClearFormats() afterwards drops the temporary text-format style so the cell
# keeps its original (default) style, matching the source workbook.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "61.564.73"
$ws.Range("E2").Value = "  -3.47%  "
$ws.Range("D3").Value = "3.002.21"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "537.58"
$ws.Range("E5").Value = "  -0.57%  "
Set-TextValue $ws.Range("D6") "134.45"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "2.995.48"
$ws.Range("E8").Value = "  -2.94%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -5.17%  "
$ws.Range("E11").Value = "  -3.00%  "
Set-TextValue $ws.Range("D12") "0.447"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("E13").Value = "  -2.34%  "
Set-TextValue $ws.Range("D14") "33.93"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "3.486.53"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "61.571.80"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "3.000.63"
$ws.Range("E18").Value = "  -2.95%  "
Set-TextValue $ws.Range("D19") "6.62"
$ws.Range("E19").Value = "  -1.39%  "
Set-TextValue $ws.Range("D20") "467.66"
$ws.Range("E20").Value = "  -4.50%  "
$ws.Range("E21").Value = "  -2.14%  "
Set-TextValue $ws.Range("D22") "0.675"
$ws.Range("E22").Value = "  -3.93%  "
Set-TextValue $ws.Range("D23") "6.93"
$ws.Range("E23").Value = "  -3.69%  "
Set-TextValue $ws.Range("D24") "80.36"
$ws.Range("E24").Value = "  +0.58%  "
Set-TextValue $ws.Range("D25") "11.97"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E28").Value = "  -6.70%  "
Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("B30").Value = "Mantle"
$ws.Range("C30").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D30") "1.15"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "1.88"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D32") "25.64"
$ws.Range("E32").Value = "  -2.58%  "
Set-TextValue $ws.Range("D33") "5.49"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("E34").Value = "  -5.44%  "
Set-TextValue $ws.Range("D35") "54.93"
$ws.Range("E35").Value = "  -3.45%  "
Set-TextValue $ws.Range("D36") "5.90"
$ws.Range("E36").Value = "  -3.00%  "
Set-TextValue $ws.Range("D37") "454.62"
$ws.Range("E37").Value = "  -8.11%  "
$ws.Range("D38").Value = "3.170.92"
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("E39").Value = "  +2.11%  "
Set-TextValue $ws.Range("D40") "0.0786"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("E41").Value = "  -3.71%  "
Set-TextValue $ws.Range("D42") "8.12"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("E43").Value = "  -7.40%  "
Set-TextValue $ws.Range("D44") "26.84"
$ws.Range("E44").Value = "  +7.02%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -5.82%  "
$ws.Range("E47").Value = "  -4.04%  "
Set-TextValue $ws.Range("D48") "119.06"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("E50").Value = "  -7.91%  "
$ws.Range("E51").Value = "  +5.95%  "
